$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 109, shifting existing rows 109:197 down to 110:198
$ws.Rows.Item(109).Insert()

# Populate the newly inserted row 109 with the new record's values
$ws.Range("A109").Value = 10
$ws.Range("B109").Value = "Vega Modelo de Temuco"
$ws.Range("C109").Value = "La Araucanía"
$ws.Range("D109").Value = 44827
$ws.Range("E109").Value = 9
$ws.Range("F109").Value = 100112012
$ws.Range("G109").Value = "Espinaca"
$ws.Range("H109").Value = "Sin especificar"
$ws.Range("I109").Value = "Primera"
$ws.Range("J109").Value = 30
$ws.Range("K109").Value = 10000
$ws.Range("L109").Value = 10000
$ws.Range("M109").Value = 10000
$ws.Range("N109").Value = "$/docena de atados"
$ws.Range("O109").Value = "Región de La Araucanía"
$ws.Range("P109").Value = 3333
$ws.Range("Q109").Value = 3
$ws.Range("R109").Value = "Hortaliza"
